$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("table_definitions")
Write-Host $ws.Name
